$p = $ppt.ActivePresentation

# Slide 1 title: "Example numbering MWE"
# Split the run "Example " into "Example" + " ", and "numbering " into
# "numbering" + " ", leaving "MWE" as its own run (matches the reverted
# "consolidate text run nodes" commit: each word/space becomes its own run).
$s1 = $p.Slides.Item(1)
$sh1 = $s1.Shapes.Item(1)
$tr1 = $sh1.TextFrame.TextRange
$tr1.Characters(1, 7).Text = "Example"
$tr1.Characters(9, 9).Text = "numbering"

# Slide 2 title: "A second slide"
# Split the run "A " into "A" + " ", and "second " into "second" + " ",
# leaving "slide" as its own run.
$s2 = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(1)
$tr2 = $sh2.TextFrame.TextRange
$tr2.Characters(1, 1).Text = "A"
$tr2.Characters(3, 6).Text = "second"
